$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Coluna 1"
$ws.Range("B1").Value = "Coluna 2"

$ws.Range("A2").Select()
